$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 14:20"

# Full refreshed table: row, Ciudad, Casos totales, Casos activos, Recuperados, Muertes
$rows = @(
    @(4, 'Madrid', 32155, 12400, 15580, 4175),
    @(5, 'Cataluña', 19991, 5701, 12441, 1849),
    @(6, 'Bizkaia/Vizcaya', 3382, 2165, 2289, 168),
    @(7, 'Valencia/Valencia', 3291, 232, 2849, 210),
    @(8, 'Ciudad Real', 2807, 492, 2340, 272),
    @(9, 'Castilla-La Mancha', 2780, 71, 2446, 263),
    @(10, 'Navarra', 2682, 278, 2263, 141),
    @(11, 'Alacant/Alicante', 2320, 172, 1960, 188),
    @(12, 'Araba/Alava', 2250, 2165, 1376, 149),
    @(13, 'Albacete', 2098, 492, 1806, 170),
    @(14, 'La Rioja', 2083, 641, 1341, 101),
    @(15, 'Zaragoza', 2068, 259, 1654, 155),
    @(16, 'A Coruña', 1969, 333, 1788, 67),
    @(17, 'Toledo', 1673, 492, 1339, 234),
    @(18, 'Malaga', 1644, 93, 1458, 93),
    @(19, 'Pontevedra', 1536, 333, 1411, 30),
    @(20, 'Tenerife', 1444, 77, 1241, 68),
    @(21, 'Salamanca', 1413, 272, 986, 155),
    @(22, 'Asturias', 1384, 135, 1180, 69),
    @(23, 'Sevilla', 1371, 20, 1294, 57),
    @(24, 'Cantabria', 1268, 60, 1148, 60),
    @(25, 'Granada', 1230, 15, 1129, 86),
    @(26, 'Caceres', 1212, 45, 1012, 155),
    @(27, 'Gipuzkoa/Guipuzcoa', 1206, 2165, 639, 52),
    @(28, 'Valladolid', 1188, 333, 750, 105),
    @(29, 'Murcia', 1084, 45, 997, 42),
    @(30, 'Leon', 1066, 276, 663, 127),
    @(31, 'Aragon', 907, 29, 838, 40),
    @(32, 'Burgos', 895, 267, 542, 86),
    @(33, 'Segovia', 886, 235, 563, 88),
    @(34, 'Guadalajara', 796, 492, 644, 105),
    @(35, 'Jaen', 788, 17, 732, 39),
    @(36, 'Cordoba', 782, 4, 753, 25),
    @(37, 'Ourense', 751, 333, 660, 22),
    @(38, 'Cadiz', 697, 17, 661, 19),
    @(39, 'Castello/Castellon', 697, 28, 624, 45),
    @(40, 'Soria', 686, 95, 545, 46),
    @(41, 'Badajoz', 625, 94, 505, 26),
    @(42, 'Lugo', 586, 333, 520, 11),
    @(43, 'Avila', 560, 151, 345, 64),
    @(44, 'Gran Canaria', 384, 77, 1241, 354),
    @(45, 'Palencia', 383, 59, 300, 24),
    @(46, 'Huesca', 349, 35, 296, 18),
    @(47, 'Cuenca', 308, 492, 207, 73),
    @(48, 'Almeria', 290, 14, 258, 18),
    @(49, 'Teruel', 283, 26, 232, 25),
    @(50, 'Zamora', 278, 61, 189, 28),
    @(51, 'Huelva', 224, 2, 216, 6),
    @(52, 'Mallorca', 210, 18, 194, 12),
    @(53, 'La Palma', 68, 77, 1241, 56),
    @(54, 'Melilla', 62, 0, 61, 1),
    @(55, 'Lanzarote', 59, 77, 1241, 46),
    @(56, 'Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena', 58, 0, 58, 3),
    @(57, 'Ceuta', 51, 0, 50, 1),
    @(58, 'Fuerteventura', 32, 77, 1241, 32),
    @(59, 'Ibiza', 21, 18, 20, 1),
    @(60, 'Menorca', 15, 18, 13, 0),
    @(61, 'La Gomera', 8, 77, 1241, 1),
    @(62, 'Arroyo de la Luz', 7, 0, 7, 0),
    @(63, 'El Hierro', 3, 77, 1241, 3),
    @(64, 'Formentera', 0, 10, 0, 8)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
